$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Set B3 (TransponderReceiver column, row with Step=2) to "S"
$ws.Range("B3").Value = "S"

# Update the active selection to B7
$ws.Range("B7").Select()
